# Updates the cryptos list (Price / Volume(1h) columns, plus a handful of
# swapped coin rows) to reflect the latest scrape, matching the commit
# "Updated cryptos list on Sun Oct 22 16:50:02 UTC 2023 with GitHub Actions".
#
# Price values are stored as text in this sheet (e.g. "29.867.43",
# "1.00"), so for any new price string that Excel would otherwise parse
# as a plain number we force the cell to Text format first so the COM
# layer keeps it as a literal string instead of silently converting it
# to a numeric value (which would also drop meaningful trailing zeros).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.863.40'
$ws.Range('E2').Value = '  +0.64%  '
$ws.Range('D3').Value = '1.630.30'
$ws.Range('E3').Value = '  +0.91%  '
$ws.Range('E4').Value = '  +0.69%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.72'
$ws.Range('E5').Value = '  +1.01%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.521'
$ws.Range('E6').Value = '  -0.16%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.72%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '28.74'
$ws.Range('E8').Value = '  -0.87%  '
$ws.Range('E9').Value = '  -0.10%  '
$ws.Range('E10').Value = '  -0.18%  '
$ws.Range('E11').Value = '  -1.16%  '
$ws.Range('D12').Value = '1.864.38'
$ws.Range('E12').Value = '  +0.89%  '
$ws.Range('D13').Value = '1.635.49'
$ws.Range('E13').Value = '  +1.35%  '
$ws.Range('E14').Value = '  -0.31%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '9.42'
$ws.Range('E15').Value = '  +4.89%  '
$ws.Range('D16').Value = '29.884.09'
$ws.Range('E16').Value = '  +0.69%  '
$ws.Range('E17').Value = '  -0.83%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '65.43'
$ws.Range('E18').Value = '  +1.89%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '240.65'
$ws.Range('E19').Value = '  -0.66%  '
$ws.Range('D20').Value = '0.0₃0702'
$ws.Range('E20').Value = '  -1.41%  '
$ws.Range('E21').Value = '  +0.54%  '
$ws.Range('E22').Value = '  +1.05%  '
$ws.Range('E23').Value = '  +0.53%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.17'
$ws.Range('E24').Value = '  +2.60%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '157.71'
$ws.Range('E25').Value = '  +0.75%  '
$ws.Range('E26').Value = '  -1.04%  '
$ws.Range('E27').Value = '  -0.92%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.59'
$ws.Range('E28').Value = '  -0.15%  '
$ws.Range('E29').Value = '  +0.60%  '
$ws.Range('E30').Value = '  +0.24%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.10'
$ws.Range('E31').Value = '  +2.63%  '
$ws.Range('E32').Value = '  +1.74%  '
$ws.Range('E33').Value = '  -1.21%  '
$ws.Range('D34').Value = '1.426.61'
$ws.Range('E34').Value = '  -0.67%  '
$ws.Range('E35').Value = '  +3.15%  '
$ws.Range('E36').Value = '  -2.79%  '
$ws.Range('E37').Value = '  -3.84%  '
$ws.Range('E38').Value = '  +0.08%  '
$ws.Range('E39').Value = '  +0.16%  '
$ws.Range('B40').Value = 'ImmutableX'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.556'
$ws.Range('E40').Value = '  -0.19%  '
$ws.Range('B41').Value = 'Aave'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '74.91'
$ws.Range('E41').Value = '  +7.05%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.832'
$ws.Range('E42').Value = '  +0.25%  '
$ws.Range('B43').Value = 'Kaspa'
$ws.Range('C43').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0499'
$ws.Range('E43').Value = '  -0.48%  '
$ws.Range('B44').Value = 'RenderToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.98'
$ws.Range('E44').Value = '  +0.59%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.00'
$ws.Range('E45').Value = '  +0.71%  '
$ws.Range('E46').Value = '  +0.59%  '
$ws.Range('B47').Value = 'RocketPoolETH'
$ws.Range('C47').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D47').Value = '1.771.81'
$ws.Range('E47').Value = '  +0.84%  '
$ws.Range('B48').Value = 'FraxShare'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '5.33'
$ws.Range('E48').Value = '  -2.26%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '48.97'
$ws.Range('E49').Value = '  -8.59%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '92.30'
$ws.Range('E51').Value = '  +4.23%  '
